# Add a new "Date" column at the front of the sheet (Sheet1), shifting the
# existing Account/Campaign/Keyword/Match Type/Cost/Clicks/Conversions/
# Impressions/Impression Share columns one to the right (B:J), and leave the
# active selection on C5 (matches the saved workbook view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A; this shifts A:I -> B:J.
$ws.Range("A1").EntireColumn.Insert()

# Populate the new header cell.
$ws.Range("A1").Value = "Date"

# Restore the saved cursor/selection position.
$ws.Range("C5").Select() | Out-Null
